$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.047.92"
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("D3").Value = "2.615.03"
$ws.Range("E3").Value = "  +3.44%  "
$ws.Range("D5").Value = "605.25"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").Value = "178.93"
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.524"
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("D9").Value = "2.612.17"
$ws.Range("E9").Value = "  +3.36%  "
$ws.Range("E10").Value = "  +12.59%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "0.347"
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("D13").Value = "5.05"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").Value = "3.139.64"
$ws.Range("E14").Value = "  +6.38%  "
$ws.Range("D15").Value = "0.0000185"
$ws.Range("E15").Value = "  +6.93%  "
$ws.Range("D16").Value = "26.65"
$ws.Range("E16").Value = "  +1.56%  "
$ws.Range("D17").Value = "70.991.51"
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("D18").Value = "2.631.92"
$ws.Range("E18").Value = "  +4.04%  "
$ws.Range("D19").Value = "380.25"
$ws.Range("E19").Value = "  +5.60%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "7.88"
$ws.Range("E20").Value = "  +4.53%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "11.45"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "1.98"
$ws.Range("E23").Value = "  +16.70%  "
$ws.Range("D24").Value = "72.25"
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("D25").Value = "4.42"
$ws.Range("E25").Value = "  +4.67%  "
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  +10.85%  "
$ws.Range("D28").Value = "2.752.06"
$ws.Range("E28").Value = "  +3.66%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "546.68"
$ws.Range("E30").Value = "  +5.09%  "
$ws.Range("D31").Value = "0.0₃0953"
$ws.Range("E31").Value = "  +6.49%  "
$ws.Range("D32").Value = "8.04"
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("D33").Value = "1.33"
$ws.Range("E33").Value = "  +6.76%  "
$ws.Range("E34").Value = "  +2.48%  "
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "166.36"
$ws.Range("E36").Value = "  +2.34%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "19.17"
$ws.Range("E37").Value = "  +3.71%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "0.114"
$ws.Range("E38").Value = "  -5.13%  "
$ws.Range("E39").Value = "  +2.43%  "
$ws.Range("D40").Value = "1.40"
$ws.Range("E40").Value = "  +6.25%  "
$ws.Range("D41").Value = "1.87"
$ws.Range("E41").Value = "  +5.42%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "2.61"
$ws.Range("E43").Value = "  +8.77%  "
$ws.Range("D44").Value = "5.04"
$ws.Range("E44").Value = "  +3.99%  "
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("D46").Value = "39.97"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("D47").Value = "152.64"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("E50").Value = "  +4.89%  "
$ws.Range("D51").Value = "0.0⁦0262"
$ws.Range("E51").Value = "  +3.05%  "
